$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "操作命令表" (sheet1): add 3 new rows documenting new ruler commands,
# and refine the wording of the existing "UpdateRuler" row.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("操作命令表")

# Row 23 (existing "UpdateRuler" entry): clarify the description text.
$ws1.Range("G23").Value = "更新标尺数据，或者是否上下行分设（但名称不变）"

# Row 24: "ChangeRulerName" command.
$ws1.Range("A24").Value = 22
$ws1.Range("B24").Value = "更改标尺名称"
$ws1.Range("C24").Value = "ChangeRulerName"
$ws1.Range("D24").Value = "支持"
$ws1.Range("E24").Value = "否"
$ws1.Range("G24").Value = "修改标尺的名称（但数据不变）"

# Row 25: "RemoveRuler" command.
$ws1.Range("A25").Value = 23
$ws1.Range("B25").Value = "删除标尺"
$ws1.Range("C25").Value = "RemoveRuler"
$ws1.Range("D25").Value = "支持"
$ws1.Range("E25").Value = "否"
$ws1.Range("G25").Value = "删除指定标尺，如果是排图标尺，同时把该线路排图标尺设为空"

# Row 26: "AddNewRuler" command.
$ws1.Range("A26").Value = 24
$ws1.Range("B26").Value = "新建空白标尺"
$ws1.Range("C26").Value = "AddNewRuler"
$ws1.Range("D26").Value = "支持"
$ws1.Range("E26").Value = "否"
$d26 = Get-Date -Year 2021 -Month 8 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws1.Range("F26").Value = $d26.Date
$ws1.Range("G26").Value = "添加空白标尺，同时打开编辑面板"

# Bring sheet1 to the front and update selection / scroll position.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("G26").Select()

# ---------------------------------------------------------------------------
# Sheet "进度" (sheet3): move the viewport / selection further down, since the
# content referenced there just shifted while no cell content actually
# changes on this sheet.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("进度")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("C33").Select()

# Re-activate the command table sheet, which ends up the active tab.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("G26").Select()
